$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update the "Förändrad" (changed) date column C for rows 2-7
# from serial date 45183 (2023-09-14) to 45184 (2023-09-15)
foreach ($row in 2..7) {
    $ws.Cells.Item($row, 3).Value = 45184
}
